$d = $word.ActiveDocument

# Replace all occurrences of "June 04, 2022" with "June 05, 2022"
$d.Content.Find.Execute("June 04, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "June 05, 2022", 2)

# Replace "August 03, 2022" with "August 04, 2022"
$d.Content.Find.Execute("August 03, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "August 04, 2022", 2)
